$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-number / update scenario rows (B:C), inserting extra scenarios
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "OK"

$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "FALLA"

$ws.Range("B5").Value = 6
$ws.Range("C5").Value = "OK"

$ws.Range("B6").Value = 7
$ws.Range("C6").Value = "OK"

$ws.Range("B7").Value = 8
$ws.Range("C7").Value = "OK"

$ws.Range("B8").Value = 9
$ws.Range("C8").Value = "OK"

$ws.Range("B9").Value = 10
$ws.Range("C9").Value = "OK"

$ws.Range("B10").Value = 12
$ws.Range("C10").Value = "OK"

$ws.Range("B11").Value = 16
$ws.Range("C11").Value = "OK"

$ws.Range("B12").Value = 17
$ws.Range("C12").Value = "OK"

$ws.Range("B13").Value = 18
$ws.Range("C13").Value = "OK"

$ws.Range("C3").Select()
